# New submission synced into the "JSS 3E" results sheet:
#   - the previous last row (row 4, Admission No "7") is normalized from
#     text to a real number, matching how the sync script stores numeric
#     Admission Nos once cleaned.
#   - a brand new response row (row 5) is appended for ESTHER YAGA.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3E")

# Row 4: Admission No "7" was stored as text; the new sync fixes it to a number.
$ws.Cells.Item(4, 3).Value = 7

# Row 5: new form submission appended at the bottom of the sheet.
$ws.Cells.Item(5, 1).Value = "2026-02-08 17:46:17"
$ws.Cells.Item(5, 2).Value = "ESTHER YAGA"

# Admission No "32" must stay text (leading context in the source sheet is
# non-numeric for this entrant), so force text formatting before writing it.
$ws.Cells.Item(5, 3).NumberFormat = "@"
$ws.Cells.Item(5, 3).Value = "32"

$ws.Cells.Item(5, 4).Value = 8
